$d = $word.ActiveDocument

# --- Step 1: update the text of the two existing runs ---
# Original paragraph text is "Firs" + "t doc" (runs split by the _GoBack bookmark).
# Target paragraph text is "First doc" + " Hello world!!!", with the bookmark
# relocated to the very end of the paragraph (after both runs).

# Replace the second run's text ("t doc" -> " Hello world!!!") using its
# known character offsets; this keeps it as a distinct run.
$r2 = $d.Range(4, 9)
$r2.Text = " Hello world!!!"

# Replace the first run's text ("Firs" -> "First doc").
$r1 = $d.Range(0, 4)
$r1.Text = "First doc"

# --- Step 2: move the _GoBack bookmark to the end of the paragraph ---
# Adding a zero-length bookmark exactly at the start/end boundary of the
# paragraph's text collapses (engine quirk) into spanning the whole
# paragraph, and setting .Text on a zero-length range touching the very end
# also merges the surrounding runs together. To avoid both pitfalls we:
#   1. Append a sacrificial character with InsertAfter (keeps runs distinct).
#   2. Anchor the bookmark at the now-interior position just before it.
#   3. Remove the sacrificial character again.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$end = $d.Content.End
$tmp = $d.Range($end - 1, $end - 1)
$tmp.InsertAfter("X")

$newEnd = $d.Content.End
$pos = $newEnd - 2
$target = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $target)

$finalEnd = $d.Content.End
$sacrifice = $d.Range($finalEnd - 2, $finalEnd - 1)
$sacrifice.Text = ""
